$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-21 04:51:57"

$wsZhCn.Range("H3").Value = "2016-08-21 04:51:53"
$wsZhCn.Range("K3").Value = "2016-08-21 04:52:14"

$wsDeDe.Range("H3").Value = "2016-08-21 04:51:57"
$wsDeDe.Range("K3").Value = "2016-08-21 04:52:21"
